$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Junio" -> "Febrero" (keeps its own run) and move the "_GoBack"
#    bookmark from the end of the section to right after that run.
# ------------------------------------------------------------------

# Remove the existing "_GoBack" bookmark (it currently sits at the end
# of the section, right before the following paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rng = $d.Content
$found = $rng.Find.Execute("Junio")
$start = $rng.Start
# Toggle a throwaway direct-formatting flag on/off around the text
# replacement: it keeps this run distinct from its identically
# formatted neighbour while the text is rewritten, then the reset
# restores the original (shared) formatting without re-merging the
# run into its neighbours.
$rng.Bold = 1
$rng.Text = "Febrero"
$endPos = $rng.End
$d.Range($start, $endPos).Bold = 0

# Re-insert the bookmark immediately after the run that now reads
# "Febrero".
$d.Bookmarks.Add("_GoBack", $d.Range($endPos, $endPos))

# ------------------------------------------------------------------
# 2) Collapse runs that were split apart purely by spell-check
#    proofing marks (<w:proofErr .../>) back into single runs. Doing
#    a Find/Replace of the full phrase (old text == new text) causes
#    the engine to re-flow the range into one run and drop the now
#    orphaned proofErr markers, matching the target XML.
# ------------------------------------------------------------------

$mergePhrases = @(
    "Hay una acción correctiva que impacte los requerimientos o necesidades de información de los stakeholders.",
    "El project manager realiza la invitación.",
    "La invitación se envía a los interesados de la reunión, utilizando el medio (correo electrónico, whatsapp).",
    "Procedimiento para él envió de información por whatsapp:",
    "el Project manager o el equipo de trabajo envía información o documentación exclusivamente del proyecto al grupo de whatsapp.",
    "En casos particulares deberán tratarse en el whatsapp personal",
    "BBB = Abreviatura del Tipo de Documento= pch, sst, wbs, dwbs, org, ram, etc. ",
    "DDDD = Formato del Archivo=doc, exe, pdf, mpp, etc."
)

foreach ($phrase in $mergePhrases) {
    $r = $d.Content
    $r.Find.Execute($phrase, $true, $false, $false, $false, $false, $true, 0, $false, $phrase, 1) | Out-Null
}
